# "Generate Report for Handback" — populate the per-locale handback columns
# (Latest Target File / Latest Handback File / Latest Handback DateTime) and
# flip the Status column from "Ready for handoff" to the handed-back state.

$wb = $excel.ActiveWorkbook

function Update-LocaleSheet {
    param([string]$SheetName, [string]$HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # ---- 1. Status column: report now reflects the handback ----
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # ---- 2. Latest Handback DateTime (column H) ----
    $ws.Range("H2").Value = $HandbackDateTime
    $ws.Range("H3").Value = $HandbackDateTime

    # ---- 3. Remember current hyperlinks (ref -> address/display) so we can
    #         rebuild the collection in a stable, row-major order after
    #         adding the two new ones per row. ----
    $hlInfo = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $ref = $hl.Range.Address($false, $false)
        $hlInfo[$ref] = @($hl.Address, $hl.TextToDisplay)
    }

    $mdAddress2  = $hlInfo["A2"][0]
    $mdDisplay2  = $hlInfo["A2"][1]
    $xlfAddress2 = $hlInfo["D2"][0]
    $xlfDisplay2 = $hlInfo["D2"][1]

    $mdAddress3  = $hlInfo["A3"][0]
    $mdDisplay3  = $hlInfo["A3"][1]
    $xlfAddress3 = $hlInfo["D3"][0]
    $xlfDisplay3 = $hlInfo["D3"][1]

    # ---- 4. Latest Target File (F) / Latest Handback File (G) values ----
    # (Hyperlinks.Add below both writes the display text into the cell AND
    # applies the workbook's hyperlink formatting, same as the A/B/D cells.)
    $hlInfo["F2"] = @($mdAddress2, $mdDisplay2)
    $hlInfo["G2"] = @($xlfAddress2, $xlfDisplay2)
    $hlInfo["F3"] = @($mdAddress3, $mdDisplay3)
    $hlInfo["G3"] = @($xlfAddress3, $xlfDisplay3)

    # ---- 5. Rebuild hyperlinks in row-major order: A2,B2,D2,F2,G2,A3,B3,D3,F3,G3 ----
    $ws.Hyperlinks.Delete()
    $order = @("A2", "B2", "D2", "F2", "G2", "A3", "B3", "D3", "F3", "G3")
    foreach ($addr in $order) {
        $pair = $hlInfo[$addr]
        $ws.Hyperlinks.Add($ws.Range($addr), $pair[0], "", "", $pair[1]) | Out-Null
    }
}

Update-LocaleSheet "zh-cn" "2016-03-12 00:43:25"
Update-LocaleSheet "de-de" "2016-03-12 00:43:30"
